$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Significant Components" — reorder the variable names inside the
# bracketed lists for factors 1, 2, 3 and 5 (factor 4 / C5 is untouched).
# ---------------------------------------------------------------------------
$wsSig = $wb.Worksheets.Item("Significant Components")
$wsSig.Range("C2").Value = "['QHISPC' 'QNOHLTH' 'QESL' 'QEXTRCT' 'QEDLESHI' 'PPUNIT' 'QPOVTY' 'PERCAP']"
$wsSig.Range("C3").Value = "['QPOVTY' 'QFAM' 'QSERV' 'QFHH' 'QBLACK' 'QRICH' 'PERCAP']"
$wsSig.Range("C4").Value = "['QFEMALE' 'QFEMLBR' 'QAGEDEP']"
$wsSig.Range("C6").Value = "['PPUNIT' 'QNOAUTO' 'QRENTER']"

# ---------------------------------------------------------------------------
# Sheet: "Loading Factors" — rows re-sorted (variables re-ranked by loading)
# and every numeric loading refreshed with the re-computed values.
# ---------------------------------------------------------------------------
$wsLoad = $wb.Worksheets.Item("Loading Factors")

$loadRows = @(
  @("QHISPC",   0.8803370950973293,  0.1202941227691,     -0.1034740134268362,   -0.2055870801357961,  -0.1433395955485494),
  @("QNOHLTH",  0.742742621222788,   0.3825258295162977,  -0.09377154927699263,  -0.1242320407102496,   0.1047037166043549),
  @("QESL",     0.8669870467330312,  0.1055887583105761,  -0.1218582076081707,   -0.164540619176934,    0.1246202491228288),
  @("QEXTRCT",  0.7518642472608691,  0.1080633736403501,  -0.1841536377102106,   -0.04880562857143168,  0.02811186015142211),
  @("QEDLESHI", 0.8692554803215153,  0.2820098250435661,  -0.09611805776404492,  -0.04435359349460323,  0.0176660870593831),
  @("PPUNIT",   0.5203431762810202,  0.1723556429544132,  -0.00853118410586568,  -0.04118569700388107,  -0.6351977639248791),
  @("QPOVTY",   0.4926703477871059,  0.4936970987725287,   0.01006171317970551,  -0.139980913852795,     0.3593032030497998),
  @("QFAM",     0.1243576388719362,  0.6629400763354163,   0.0760197032434357,   -0.1184125657849894,    0.2379118995421976),
  @("QSERV",    0.3709861022228685,  0.5379072856593253,   0.005759199727788074, -0.08757769907133597,   0.155983972952657),
  @("QFHH",     0.2285252661260082,  0.7131346559507746,   0.2214646558077026,   -0.07783278425410714,  -0.04601151511459178),
  @("QBLACK",  -0.2749361517979486,  0.7046485575443674,   0.05019711505949858,   0.1307361815087761,    0.1670202073710436),
  @("QRICH",    0.4038778838595897,  0.6463639855284391,  -0.06158261054806106,  -0.1194534729670854,    0.002273121296869304),
  @("PERCAP",   0.5037823276587938,  0.7061433370752405,  -0.04688443646465085,  -0.1072480927953808,   -0.1020625146575289),
  @("QNOAUTO",  0.1891304089574793,  0.3968119682334907,  -0.001312813755978878,  0.08291571045634666,   0.5521493354761324),
  @("QRENTER",  0.1727227389210597,  0.3584593742508944,  -0.05642809061090642,  -0.46054991919618,      0.6579699975112047),
  @("QFEMALE", -0.1201784096346595,  0.08442001695086726,  0.9476420664224633,    0.01306748968427164,  -0.03709586659921296),
  @("QFEMLBR", -0.4233895389422812,  0.2565876914832129,   0.5549136341195142,   -0.009510221564061949, -0.03382139079307862),
  @("QAGEDEP", -0.1133192017874101, -0.06739481440699845,  0.729040421503802,     0.4777003061560329,    0.06612984087372405),
  @("MEDAGE",  -0.3065470999719298, -0.3702019232495327,   0.04864970237269507,   0.6451312943906372,    0.09075703803451746),
  @("QSSBEN",  -0.08069646390827452, 0.0708979599700322,   0.1345297566914315,    0.847625101158753,    -0.1136918263240227)
)

$r = 2
foreach ($row in $loadRows) {
    $wsLoad.Cells.Item($r, 1).Value = $row[0]
    $wsLoad.Cells.Item($r, 2).Value = $row[1]
    $wsLoad.Cells.Item($r, 3).Value = $row[2]
    $wsLoad.Cells.Item($r, 4).Value = $row[3]
    $wsLoad.Cells.Item($r, 5).Value = $row[4]
    $wsLoad.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet: "All Refactor Variances" — refresh the re-computed F1 group figures
# in columns I:M (rows 2-5); the F0 group in columns B:H is unchanged.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Refactor Variances")

$allVarRows = @(
  @(4.983165037844279, 3.692480322986977, 1.902149201692064, 1.763128705847649, 1.467933700943554),
  @(0.249158251892214, 0.1846240161493488, 0.09510746008460322, 0.08815643529238243, 0.07339668504717768),
  @(0.249158251892214, 0.4337822680415628, 0.5288897281261661, 0.6170461634185485, 0.6904428484657261),
  @(0.3608673077661435, 0.2673994184451512, 0.1377484904014099, 0.1276810028350356, 0.1063037805522598)
)

$r = 2
foreach ($row in $allVarRows) {
    $wsAll.Cells.Item($r, 9).Value  = $row[0]
    $wsAll.Cells.Item($r, 10).Value = $row[1]
    $wsAll.Cells.Item($r, 11).Value = $row[2]
    $wsAll.Cells.Item($r, 12).Value = $row[3]
    $wsAll.Cells.Item($r, 13).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet: "Final Variances" — same refreshed figures, columns B:F.
# ---------------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Final Variances")

$finalVarRows = @(
  @(4.983165037844279, 3.692480322986977, 1.902149201692064, 1.763128705847649, 1.467933700943554),
  @(0.249158251892214, 0.1846240161493488, 0.09510746008460322, 0.08815643529238243, 0.07339668504717768),
  @(0.249158251892214, 0.4337822680415628, 0.5288897281261661, 0.6170461634185485, 0.6904428484657261),
  @(0.3608673077661435, 0.2673994184451512, 0.1377484904014099, 0.1276810028350356, 0.1063037805522598)
)

$r = 2
foreach ($row in $finalVarRows) {
    $wsFinal.Cells.Item($r, 2).Value = $row[0]
    $wsFinal.Cells.Item($r, 3).Value = $row[1]
    $wsFinal.Cells.Item($r, 4).Value = $row[2]
    $wsFinal.Cells.Item($r, 5).Value = $row[3]
    $wsFinal.Cells.Item($r, 6).Value = $row[4]
    $r++
}

# ---------------------------------------------------------------------------
# Sheet: "Included and Excluded" — the "include" list is re-ordered to match
# the new variable ranking; "exclude" list is unchanged.
# ---------------------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Included and Excluded")
$wsInc.Range("B2").Value = "[['QHISPC', 'QNOHLTH', 'QESL', 'QEXTRCT', 'QEDLESHI', 'PPUNIT', 'QPOVTY', 'PERCAP', 'QFAM', 'QSERV', 'QFHH', 'QBLACK', 'QRICH', 'QFEMALE', 'QFEMLBR', 'QAGEDEP', 'QRENTER', 'MEDAGE', 'QSSBEN', 'QNOAUTO']]"

Write-Host "edit complete"
